$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.029143929481506
$ws.Range("B1").Value = 2.122998714447021
$ws.Range("C1").Value = 3.96337628364563
$ws.Range("D1").Value = 0.9487218856811523
$ws.Range("E1").Value = 0.8074910640716553
